$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.027.41"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "3.205.66"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.205.11"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.507"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "3.734.47"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "66.189.45"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "3.216.54"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.111"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.38%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0909"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "484.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("E42").Value = "  +3.88%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("D45").Value = "2.953.20"
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("D46").Value = "0.0₃0643"
$ws.Range("E46").Value = "  +3.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").Value = "  +2.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.23%  "
